# Rename the worksheet "Property1" -> "DataNode" to unify the
# DataNode/DataTable/Entity naming convention described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "DataNode"

# Update the current selection to match the saved view state: the
# active cell stays A9 (top-left of the frozen pane) but the selected
# range now spans A9:N35.
$ws.Range("A9:N35").Select() | Out-Null
